# Apply the arithmetic-problem text replacements described in the diff.
# Each old value is unique in the document and each new value is unique
# too, so a straightforward ordered sequence of single Find/Replace
# operations (wrap = 0 / find whole story, MatchCase = true) is safe.

$d = $word.ActiveDocument

$replacements = @(
    @("798×4=3192", "572×3=1716"),
    @("895×8=7160", "795×8=6360"),
    @("979×9=8811", "542×6=3252"),
    @("460×6=2760", "739×7=5173"),
    @("191×8=1528", "980×4=3920"),
    @("397×2=794",  "781×4=3124"),
    @("108×2=216",  "177×3=531"),
    @("321×9=2889", "195×3=585"),
    @("981×4=3924", "111×2=222"),
    @("527×5=2635", "388×6=2328"),
    @("252×6=1512", "578×2=1156"),
    @("976×5=4880", "887×7=6209"),
    @("133×8=1064", "937×9=8433"),
    @("523×8=4184", "582×6=3492"),
    @("125×7=875",  "136×9=1224"),
    @("794×6=4764", "985×3=2955"),
    @("133×5=665",  "590×4=2360"),
    @("279×5=1395", "333×7=2331"),
    @("171×2=342",  "108×4=432"),
    @("678×6=4068", "529×2=1058"),
    @("457×5=2285", "101×7=707"),
    @("670×2=1340", "232×9=2088"),
    @("250×6=1500", "814×9=7326"),
    @("124×6=744",  "806×6=4836"),
    @("474×2=948",  "270×2=540")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)

    if (-not $found) {
        Write-Host "WARNING: could not find '$old'"
    }
}

$d.Save()
